# Added Historical VaR as a second step
# Update computed VaR values across the three report sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: Total_VaR ---
$wsTotal = $wb.Worksheets.Item("Total_VaR")
$wsTotal.Range("A2").Value = 1524154.456000008

# --- Sheet: VaR_by_BUSINESS_LINE ---
$wsBL = $wb.Worksheets.Item("VaR_by_BUSINESS_LINE")
$wsBL.Range("B2").Value = 1364332.8
$wsBL.Range("B3").Value = 1897405.728000013
$wsBL.Range("B4").Value = 107847.9999999999

# --- Sheet: VaR_by_METAL ---
$wsMetal = $wb.Worksheets.Item("VaR_by_METAL")
$wsMetal.Range("B2").Value = 1490154.456000008
$wsMetal.Range("B3").Value = 107847.9999999999
